$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOLDE")

$c6 = $ws.Range("C6").Value2
$c7 = $ws.Range("C7").Value2

$ws.Range("C6").Value = "V 0/3 - " + $c6
$ws.Range("C7").Value = "V 0/3 - " + $c7

$ws.Range("C7").Select()
